$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H, matching the style used by the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New "Save" column values (rows 2-11)
$values = @(0, 1, 1, 1, 0, 1, 0, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
